# EIA Table 1.1.A monthly update: October 2016 -> November 2016 vintage.
# Adds the "November" monthly row to the Year-2016 block, and refreshes the
# "Year to Date" and "Rolling 12 Months" summary rows/labels accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Title banner: "...2006-October 2016" -> "...2006-November 2016"
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Table 1.1.A. Net Generation from Renewable Sources:  Total (All Sectors), 2006-November 2016"

# ---------------------------------------------------------------------------
# 2) Insert a new row 53 for "November" underneath the existing October row,
#    copying the number formatting from a sibling month row (January, row 43)
#    so no new cell-style entries get minted, then fill in its values.
# ---------------------------------------------------------------------------
$ws.Rows("53:53").Insert()
$ws.Range("A43:N43").Copy()
$ws.Range("A53:N53").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A53").Value = "November"
$novVals = @(19334,2458,184,3257,927,600,250,1507,18815,47331,1307,3766,3950)
for ($i = 0; $i -lt $novVals.Length; $i++) {
    $ws.Cells.Item(53, 2 + $i).Value = $novVals[$i]
}

# ---------------------------------------------------------------------------
# 3) "Year to Date" block (rows shifted down by the insert: now 54-57).
#    Row 54 is the section header (unchanged text), rows 55-57 hold the
#    refreshed 2014 / 2015 / 2016 year-to-date totals.
# ---------------------------------------------------------------------------
$ytd2014 = @(166944,14313,2346,38603,10330,6619,2934,14502,237038,493629,10467,24780,27126)
$ytd2015 = @(170620,20222,3101,38341,10262,6569,2902,14541,225915,492474,13225,33447,36548)
$ytd2016 = @(203453,30539,3293,36842,10551,6724,2866,15797,243220,553284,18281,48820,52113)

for ($i = 0; $i -lt $ytd2014.Length; $i++) {
    $ws.Cells.Item(55, 2 + $i).Value = $ytd2014[$i]
}
for ($i = 0; $i -lt $ytd2015.Length; $i++) {
    $ws.Cells.Item(56, 2 + $i).Value = $ytd2015[$i]
}
for ($i = 0; $i -lt $ytd2016.Length; $i++) {
    $ws.Cells.Item(57, 2 + $i).Value = $ytd2016[$i]
}

# ---------------------------------------------------------------------------
# 4) "Rolling 12 Months Ending in October" -> "...in November" (now row 58),
#    with refreshed 2015 / 2016 rolling totals in rows 59-60.
# ---------------------------------------------------------------------------
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

$roll2015 = @(185331,21158,3197,42078,11152,7179,3170,15916,248243,537424,13991,35150,38346)
$roll2016 = @(223552,31983,3419,40429,11579,7366,3164,17174,266386,605052,19195,51179,54597)

for ($i = 0; $i -lt $roll2015.Length; $i++) {
    $ws.Cells.Item(59, 2 + $i).Value = $roll2015[$i]
}
for ($i = 0; $i -lt $roll2016.Length; $i++) {
    $ws.Cells.Item(60, 2 + $i).Value = $roll2016[$i]
}
